{"js": "// Update the two-digit x two-digit multiplication answers in the table.\n// Each old \"A\u00d7B=C\" string is unique in the document, so a literal,\n// case-sensitive search safely targets exactly the run that needs updating.\nconst replacements = [\n  [\"85\u00d713=1105\", \"82\u00d724=1968\"],\n  [\"95\u00d718=1710\", \"21\u00d751=1071\"],\n  [\"83\u00d758=4814\", \"33\u00d743=1419\"],\n  [\"74\u00d757=4218\", \"97\u00d737=3589\"],\n  [\"91\u00d782=7462\", \"48\u00d794=4512\"],\n  [\"42\u00d751=2142\", \"45\u00d778=3510\"],\n  [\"40\u00d792=3680\", \"33\u00d740=1320\"],\n  [\"86\u00d796=8256\", \"99\u00d778=7722\"],\n  [\"23\u00d760=1380\", \"62\u00d747=2914\"],\n  [\"12\u00d755=660\", \"75\u00d762=4650\"],\n  [\"83\u00d765=5395\", \"99\u00d764=6336\"],\n  [\"40\u00d788=3520\", \"45\u00d718=810\"],\n  [\"72\u00d750=3600\", \"83\u00d768=5644\"],\n  [\"87\u00d773=6351\", \"57\u00d755=3135\"],\n  [\"34\u00d762=2108\", \"54\u00d755=2970\"],\n  [\"22\u00d797=2134\", \"75\u00d799=7425\"],\n  [\"69\u00d787=6003\", \"42\u00d711=462\"],\n  [\"81\u00d769=5589\", \"81\u00d744=3564\"],\n  [\"66\u00d732=2112\", \"88\u00d716=1408\"],\n  [\"58\u00d729=1682\", \"57\u00d741=2337\"],\n  [\"67\u00d772=4824\", \"26\u00d756=1456\"],\n  [\"36\u00d711=396\", \"99\u00d745=4455\"],\n  [\"59\u00d725=1475\", \"91\u00d772=6552\"],\n  [\"76\u00d719=1444\", \"97\u00d730=2910\"],\n  [\"96\u00d758=5568\", \"54\u00d741=2214\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit x two-digit multiplication answers in the table.\n# Each old \"A\u00d7B=C\" string is unique in the document, so Find/Replace\n# safely targets exactly the text that needs updating, leaving run\n# formatting (rFonts/sz) untouched.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"85\u00d713=1105\", \"82\u00d724=1968\"),\n    @(\"95\u00d718=1710\", \"21\u00d751=1071\"),\n    @(\"83\u00d758=4814\", \"33\u00d743=1419\"),\n    @(\"74\u00d757=4218\", \"97\u00d737=3589\"),\n    @(\"91\u00d782=7462\", \"48\u00d794=4512\"),\n    @(\"42\u00d751=2142\", \"45\u00d778=3510\"),\n    @(\"40\u00d792=3680\", \"33\u00d740=1320\"),\n    @(\"86\u00d796=8256\", \"99\u00d778=7722\"),\n    @(\"23\u00d760=1380\", \"62\u00d747=2914\"),\n    @(\"12\u00d755=660\",  \"75\u00d762=4650\"),\n    @(\"83\u00d765=5395\", \"99\u00d764=6336\"),\n    @(\"40\u00d788=3520\", \"45\u00d718=810\"),\n    @(\"72\u00d750=3600\", \"83\u00d768=5644\"),\n    @(\"87\u00d773=6351\", \"57\u00d755=3135\"),\n    @(\"34\u00d762=2108\", \"54\u00d755=2970\"),\n    @(\"22\u00d797=2134\", \"75\u00d799=7425\"),\n    @(\"69\u00d787=6003\", \"42\u00d711=462\"),\n    @(\"81\u00d769=5589\", \"81\u00d744=3564\"),\n    @(\"66\u00d732=2112\", \"88\u00d716=1408\"),\n    @(\"58\u00d729=1682\", \"57\u00d741=2337\"),\n    @(\"67\u00d772=4824\", \"26\u00d756=1456\"),\n    @(\"36\u00d711=396\",  \"99\u00d745=4455\"),\n    @(\"59\u00d725=1475\", \"91\u00d772=6552\"),\n    @(\"76\u00d719=1444\", \"97\u00d730=2910\"),\n    @(\"96\u00d758=5568\", \"54\u00d741=2214\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
